$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data; D and E columns are text-formatted
# to preserve exact formatting (e.g. "231.10", "  +2.45%  ").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.377.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.230.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.74%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.10"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.67%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.95"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.404"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.34%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.562.33"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.58"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.31"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.63"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.800"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.247.94"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.225.85"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0944"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.20"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.95%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.32"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "244.99"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.20%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.24%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +5.67%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.38"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.17%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.38"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.84%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.33%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.03%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.26%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.27%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.82%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0653"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.62%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.36"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.55"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.97%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.000230"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.61"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.61%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0962"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.02%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "97.08"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.53%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.461.01"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.34"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -9.02%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.16"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.62%  "

$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.75"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.26%  "

$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.08"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.49%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.15%  "
